# Generate Report for Handoff
# Update the "Latest Handoff Date/Datetime" column for rows whose status
# is "Handback transform failed" or "Ready for handoff" (i.e. rows that
# were re-stamped by the CI run that regenerated this report).

$wb = $excel.ActiveWorkbook

$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

# Overview sheet: column D holds "Latest Handoff Date"
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("D" + $r).Value = "2016-24-21 02:24:22"
}

# zh-cn sheet: column E holds "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E" + $r).Value = "2016-03-21 02:24:18"
}

# de-de sheet: column E holds "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E" + $r).Value = "2016-03-21 02:24:22"
}
